$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.328.24'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.932.12'
$ws.Range('E3').Value = '  -2.28%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'240.73"
$ws.Range('E5').Value = '  -1.56%  '
$ws.Range('D6').Value = "'0.602"
$ws.Range('E6').Value = '  -3.67%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'56.23"
$ws.Range('E8').Value = '  -4.00%  '
$ws.Range('D9').Value = "'0.357"
$ws.Range('E9').Value = '  -4.10%  '
$ws.Range('D10').Value = "'0.0830"
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('D12').Value = '2.213.52'
$ws.Range('E12').Value = '  -2.34%  '
$ws.Range('D13').Value = "'21.04"
$ws.Range('E13').Value = '  -8.45%  '
$ws.Range('D14').Value = "'0.795"
$ws.Range('E14').Value = '  -7.00%  '
$ws.Range('D15').Value = "'13.26"
$ws.Range('E15').Value = '  -4.51%  '
$ws.Range('D16').Value = "'5.09"
$ws.Range('E16').Value = '  -6.05%  '
$ws.Range('D17').Value = '1.931.35'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').Value = '36.247.13'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = '0.0₃0857'
$ws.Range('E19').Value = '  -2.82%  '
$ws.Range('D20').Value = "'68.43"
$ws.Range('E20').Value = '  -2.51%  '
$ws.Range('D21').Value = "'225.74"
$ws.Range('E21').Value = '  -3.39%  '
$ws.Range('D22').Value = "'4.92"
$ws.Range('E22').Value = '  -6.57%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = "'2.32"
$ws.Range('E24').Value = '  -7.42%  '
$ws.Range('E25').Value = '  -1.32%  '
$ws.Range('D26').Value = "'9.05"
$ws.Range('E26').Value = '  -8.23%  '
$ws.Range('D27').Value = "'160.43"
$ws.Range('E27').Value = '  -1.55%  '
$ws.Range('D28').Value = "'0.130"
$ws.Range('E28').Value = '  -1.37%  '
$ws.Range('D29').Value = "'19.02"
$ws.Range('E29').Value = '  -3.66%  '
$ws.Range('E30').Value = '  -2.58%  '
$ws.Range('E31').Value = '  -6.39%  '
$ws.Range('D32').Value = "'4.49"
$ws.Range('E32').Value = '  -7.50%  '
$ws.Range('D33').Value = "'0.0617"
$ws.Range('E33').Value = '  -9.40%  '
$ws.Range('D34').Value = "'4.12"
$ws.Range('E34').Value = '  -6.08%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  -1.38%  '
$ws.Range('D37').Value = "'5.90"
$ws.Range('E37').Value = '  -4.32%  '
$ws.Range('D38').Value = "'2.14"
$ws.Range('E38').Value = '  -4.19%  '
$ws.Range('D39').Value = "'2.97"
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D40').Value = "'0.0962"
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('D41').Value = "'2.85"
$ws.Range('E41').Value = '  -1.48%  '
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('E43').Value = '  -7.41%  '
$ws.Range('D44').Value = "'15.55"
$ws.Range('E44').Value = '  -3.10%  '
$ws.Range('D45').Value = '1.320.40'
$ws.Range('E45').Value = '  -2.85%  '
$ws.Range('D46').Value = "'1.01"
$ws.Range('E46').Value = '  -7.05%  '
$ws.Range('D47').Value = "'84.42"
$ws.Range('E47').Value = '  -7.93%  '
$ws.Range('D48').Value = "'6.97"
$ws.Range('E48').Value = '  -6.37%  '
$ws.Range('D49').Value = "'2.82"
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = '2.104.69'
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('D51').Value = "'42.93"
$ws.Range('E51').Value = '  -4.10%  '
